$d = $word.ActiveDocument

function Find-ParagraphStartingWith($doc, $marker) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($marker)) {
            return $i
        }
    }
    return -1
}

# --- Step 1: TEST30 paragraph: the Terzaghi problem description changes from
#     "2D" to "1D" now that the iterative coupled-flow/deformation solve works ---
$idx30 = Find-ParagraphStartingWith $d "TEST30"
$p30 = $d.Paragraphs.Item($idx30)
$r30 = $d.Range($p30.Range.Start, $p30.Range.End)
$r30.Find.Execute(" 2D ", $false, $false, $false, $false, $false, $true, 1, $false, " 1D ", 2) | Out-Null

# --- Step 2: remove the entire TEST31 paragraph (its description is dropped) ---
$idx31 = Find-ParagraphStartingWith $d "TEST31"
$idx33 = Find-ParagraphStartingWith $d "TEST33"
$p31 = $d.Paragraphs.Item($idx31)
$p33 = $d.Paragraphs.Item($idx33)
$delRange = $d.Range($p31.Range.Start, $p33.Range.Start)
$delRange.Delete()

# --- Step 3: relocate the "_GoBack" bookmark from its old spot inside the TEST37
#     paragraph to the start of the (now renumbered) TEST33 paragraph ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$idx33b = Find-ParagraphStartingWith $d "TEST33"
$p33b = $d.Paragraphs.Item($idx33b)
$targetPos = $p33b.Range.Start
$targetRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $targetRange)
